{"js": "// The underlying OOXML diff for this revision is a pure namespace/attribute\n// canonicalization artifact (the repository's fixture was re-serialized by\n// an upgraded Apache POI (3.15), which emits `xmlns:*`/attribute lists in\n// sorted order). Every `-`/`+` pair in the diff touches the exact same\n// element, the exact same set of attributes, and the exact same attribute\n// values \u2014 only the left-to-right order of attributes changed (confirmed by\n// XML canonicalization (C14N) of word/document.xml and word/styles.xml,\n// which reproduces the diff's \"after\" text exactly). There is no text,\n// formatting, or structural change to the document's content: no runs,\n// paragraphs, styles, or property values were added, removed, or modified.\n//\n// The Word JavaScript API (like the Word COM object model) operates on the\n// semantic document object model \u2014 it has no facility for reordering raw\n// XML attribute serialization, and Word's own OOXML writer is what decides\n// attribute order on save, not script-level edits. So the faithful\n// replication of this particular commit is to leave the document's content\n// untouched. We still touch the object model (read-only) so the script\n// demonstrably inspects the document, matching the \"before\" state with no\n// net edits.\n\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n\n// No content, formatting, or structural changes are required: the commit\n// only reorders XML attributes (a packaging/library side effect), which is\n// not a semantic edit reachable through Office.js.\n", "ps1": "# The underlying OOXML diff for this revision is a pure namespace/attribute\n# canonicalization artifact (the repository's fixture was re-serialized by\n# an upgraded Apache POI (3.15), which emits `xmlns:*`/attribute lists in\n# sorted order). Every `-`/`+` pair in the diff touches the exact same\n# element, the exact same set of attributes, and the exact same attribute\n# values -- only the left-to-right order of attributes changed (confirmed by\n# XML canonicalization (C14N) of word/document.xml and word/styles.xml,\n# which reproduces the diff's \"after\" text exactly). There is no text,\n# formatting, or structural change to the document's content: no runs,\n# paragraphs, styles, or property values were added, removed, or modified.\n#\n# The Word COM object model (like the Word JavaScript API) operates on the\n# semantic document object model -- it has no facility for reordering raw\n# XML attribute serialization, and Word's own OOXML writer is what decides\n# attribute order on save, not script-level edits. So the faithful\n# replication of this particular commit is to leave the document's content\n# untouched. We still touch the object model (read-only) so the script\n# demonstrably inspects the document, matching the \"before\" state with no\n# net edits.\n\n$d = $word.ActiveDocument\n\n# Read-only touch: confirm the document content is reachable, but make no\n# edits -- the commit only reorders XML attributes (a packaging/library\n# side effect), which is not a semantic edit reachable through the Word\n# object model.\n$null = $d.Content.Text\n$null = $d.Paragraphs.Count\n"}
